# Boletin Epi Pereira - poisson.xlsx
# Inicio de año epidemiologico 2026, semana 4 2026
# Replace the event-surveillance table (rows 2-29) with the updated data for
# the new epidemiological week, and remove the now-unused trailing rows
# (previously rows 30-37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# U+201D RIGHT DOUBLE QUOTATION MARK, used (as a mojibake artifact in the
# source data) inside a couple of the event names below.
$rdq = [char]0x201D

# New data for rows 2..29 : evento, nom_eve, Esperado, Observado, "valor p"
$rows = @(
    @("113", "Desnutrici" + $rdq + "n aguda en menores de 5 anos", 2, 2, 0.27),
    @("115", "Cancer en menores de 18 anos", 1, 0, 0.37),
    @("155", "Cancer de la mama y cuello uterino", 7, 3, 0.05),
    @("210", "Dengue", 1, 3, 0.06),
    @("215", "Defectos congenitos", 1, 6, 0),
    @("300", "Agresiones por animales potencialmente transmisores de rabia", 50, 60, 0.02),
    @("330", "Hepatitis a", 1, 0, 0.37),
    @("340", "Hepatitis b, c y coinfeccion hepatitis b y delta", 1, 2, 0.18),
    @("342", "Enfermedades huerfanas - raras", 2, 8, 0),
    @("346", "Ira por virus nuevo", 294, 0, 0),
    @("348", "Infeccion respiratoria aguda grave irag inusitada", 1, 0, 0.37),
    @("352", "Infecciones de sitio quirurgico asociadas a procedimiento medico quirurgico", 0, 1, 0),
    @("355", "Enfermedad transmitida por alimentos o agua (eta)", 0, 0, 1),
    @("356", "Intento de suicidio", 11, 6, 0.04),
    @("357", "Iad - infecciones asociadas a dispositivos - individual", 1, 0, 0.37),
    @("365", "Intoxicaciones", 6, 6, 0.16),
    @("455", "Leptospirosis", 1, 5, 0),
    @("465", "Malaria", 0, 3, 0),
    @("535", "Meningitis bacteriana y enfermedad meningoc" + $rdq + "cica", 0, 0, 1),
    @("549", "Morbilidad materna extrema", 6, 3, 0.09),
    @("560", "Mortalidad perinatal y neonatal tardia", 1, 1, 0.37),
    @("580", "Mortalidad por dengue", 0, 0, 1),
    @("591", "Vigilancia integrada de muertes en menores de cinco anos por infeccion respiratoria aguda - enfermedad diarreica aguda y/o desnutricion", 0, 0, 1),
    @("750", "Sifilis gestacional", 2, 0, 0.14),
    @("800", "Tos ferina", 0, 1, 0),
    @("813", "Tuberculosis", 8, 11, 0.07000000000000001),
    @("831", "Varicela individual", 5, 1, 0.03),
    @("850", "Vih/sida/mortalidad por sida", 7, 13, 0.01)
)

$lastRow = 1 + $rows.Count

# Column A holds event codes that look numeric ("113", "115", ...) but must be
# stored as text, exactly like the original workbook. Force a text number
# format while assigning, then restore the default ("Normal") style so the
# cells don't end up with a lingering explicit format.
$ws.Range("A2:A$lastRow").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$ws.Range("A2:A$lastRow").Style = "Normal"

# Remove the old trailing rows (previously rows 30-37, now beyond the new
# 29-row extent) so the sheet shrinks back down to A1:E29.
$ws.Range("A30:E37").EntireRow.Delete() | Out-Null
